$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.914.94"
$ws.Range("E2").Value = "  -0.67%  "
$ws.Range("D3").Value = "2.446.75"
$ws.Range("E3").Value = "  -1.61%  "
$ws.Range("E4").Value = "  -0.29%  "
$ws.Range("D5").Value = "'523.19"
$ws.Range("E5").Value = "  +0.14%  "
$ws.Range("D6").Value = "'130.77"
$ws.Range("E6").Value = "  -1.91%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "'0.563"
$ws.Range("E8").Value = "  +0.35%  "
$ws.Range("D9").Value = "2.446.18"
$ws.Range("E9").Value = "  -3.05%  "
$ws.Range("D10").Value = "'0.0981"
$ws.Range("E10").Value = "  +0.14%  "
$ws.Range("E11").Value = "  -1.79%  "
$ws.Range("E12").Value = "  -4.45%  "
$ws.Range("E13").Value = "  -2.53%  "
$ws.Range("D14").Value = "2.878.78"
$ws.Range("E14").Value = "  -3.13%  "
$ws.Range("D15").Value = "57.775.12"
$ws.Range("E15").Value = "  -0.98%  "
$ws.Range("D16").Value = "'21.68"
$ws.Range("E16").Value = "  -2.13%  "
$ws.Range("D17").Value = "'0.0000132"
$ws.Range("E17").Value = "  -1.88%  "
$ws.Range("D18").Value = "2.439.84"
$ws.Range("E18").Value = "  -3.34%  "
$ws.Range("D19").Value = "'10.26"
$ws.Range("E19").Value = "  -3.97%  "
$ws.Range("D20").Value = "'4.13"
$ws.Range("E20").Value = "  -0.68%  "
$ws.Range("D21").Value = "'315.28"
$ws.Range("E21").Value = "  -2.07%  "
$ws.Range("D22").Value = "'6.12"
$ws.Range("E22").Value = "  -0.90%  "
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("D24").Value = "'64.80"
$ws.Range("E24").Value = "  +0.23%  "
$ws.Range("E25").Value = "  -1.88%  "
$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  -0.12%  "
$ws.Range("E27").Value = "  -2.03%  "
$ws.Range("D28").Value = "'7.20"
$ws.Range("E28").Value = "  -2.65%  "
$ws.Range("D29").Value = "'174.40"
$ws.Range("E29").Value = "  +4.10%  "
$ws.Range("D30").Value = "0.0₃0734"
$ws.Range("E30").Value = "  -2.78%  "
$ws.Range("E31").Value = "  -1.93%  "
$ws.Range("D32").Value = "'6.18"
$ws.Range("E32").Value = "  -2.45%  "
$ws.Range("E33").Value = "  -5.02%  "
$ws.Range("E34").Value = "  +0.01%  "
$ws.Range("D35").Value = "'0.997"
$ws.Range("E35").Value = "  +0.40%  "
$ws.Range("D36").Value = "'17.77"
$ws.Range("E36").Value = "  -2.38%  "
$ws.Range("D37").Value = "'1.19"
$ws.Range("E37").Value = "  -5.73%  "
$ws.Range("D38").Value = "'3.77"
$ws.Range("E38").Value = "  -4.03%  "
$ws.Range("D39").Value = "'36.14"
$ws.Range("E39").Value = "  -0.81%  "
$ws.Range("E40").Value = "  -2.59%  "
$ws.Range("D41").Value = "'0.795"
$ws.Range("E41").Value = "  +2.90%  "
$ws.Range("E42").Value = "  -2.38%  "
$ws.Range("D43").Value = "'262.10"
$ws.Range("E43").Value = "  -5.23%  "
$ws.Range("E44").Value = "  -2.49%  "
$ws.Range("E45").Value = "  -3.91%  "
$ws.Range("E46").Value = "  -0.29%  "
$ws.Range("D47").Value = "'122.06"
$ws.Range("E47").Value = "  -6.45%  "
$ws.Range("D48").Value = "'0.0496"
$ws.Range("E48").Value = "  -1.23%  "
$ws.Range("D49").Value = "'0.0211"
$ws.Range("E49").Value = "  -1.41%  "
$ws.Range("E50").Value = "  -4.86%  "
$ws.Range("D51").Value = "'16.21"
$ws.Range("E51").Value = "  -4.19%  "
